$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting rows 18 and below down by one.
$ws.Rows(18).Insert()

# Populate the new row 18 with fresh data.
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 45230
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100112040
$ws.Cells.Item(18, 7).Value = "Cilantro"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 1000
$ws.Cells.Item(18, 12).Value = 1300
$ws.Cells.Item(18, 13).Value = 1180
$ws.Cells.Item(18, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(18, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 16).Value = 590
$ws.Cells.Item(18, 17).Value = 2
$ws.Cells.Item(18, 18).Value = "Hortaliza"
